# Quarterly indexing esoteric bug-fix operation
#
# Column A holds "first release" dates that should mark the day the GDP
# estimate for a given quarter was actually published (roughly six weeks
# after quarter-end), but they were mistakenly stored as the 1st of the
# release month instead of the 15th of the following month. Shift every
# date in column A (rows 2-150) forward by one month and land on the 15th,
# while leaving the value column (B) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 150; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $oldDate = $cell.Value()
    $newDate = $oldDate.AddMonths(1).AddDays(14)
    $cell.Value = $newDate
}
